$d = $word.ActiveDocument

# The last paragraph currently reads (as one run):
#   "1 GB=1000000KB ==>20971 tài liệu"
# with a zero-length "_GoBack" bookmark sitting right after that text
# (before the paragraph mark).
#
# Target layout:
#   run 1: "1 GB=1000000KB ==>20000"
#   <bookmarkStart/bookmarkEnd name="_GoBack">
#   run 2: " tài liệu"
#
# i.e. "20971" becomes "20000", and the trailing " tài liệu" is moved to a
# new run placed *after* the (still empty) _GoBack bookmark.

$bm = $d.Bookmarks("_GoBack")
$bmPos = $bm.Start

# 1) Remove the old trailing " tài liệu" (9 characters) that sits right
#    before the bookmark.
$tailRange = $d.Range($bmPos - 9, $bmPos)
$tailRange.Text = ""

# 2) Replace "20971" (the 5 characters now immediately before the bookmark)
#    with "20000".
$bm = $d.Bookmarks("_GoBack")
$bmPos = $bm.Start
$numRange = $d.Range($bmPos - 5, $bmPos)
$numRange.Text = "20000"

# 3) Insert " tài liệu" as a brand-new run right after the bookmark.
#    Toggling a character-formatting property and reverting it forces the
#    engine to keep this as its own run instead of silently re-merging it
#    with the neighbouring run that has identical formatting.
$bm = $d.Bookmarks("_GoBack")
$insPoint = $d.Range($bm.End, $bm.End)
$insPoint.InsertAfter(" tài liệu")
$insPoint.Bold = 1
$insPoint.Bold = 0

# 4) The insertion above moved the auto-tracked "_GoBack" bookmark to the
#    end of the newly inserted text; put it back where it belongs, right
#    between the two runs (i.e. 9 characters before where it is now).
$bm = $d.Bookmarks("_GoBack")
$fixPoint = $d.Range($bm.Start - 9, $bm.Start - 9)
$d.Bookmarks.Add("_GoBack", $fixPoint)
